$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.854.42'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '3.874.71'
$ws.Range('E3').Value = '  +3.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.85'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.09'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.94%  '
$ws.Range('D7').Value = '3.875.46'
$ws.Range('E7').Value = '  +3.95%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.527'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.167'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +2.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.481'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.03'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '4.497.20'
$ws.Range('E15').Value = '  +3.32%  '
$ws.Range('D16').Value = '3.857.42'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '69.891.80'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.45'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('E19').Value = '  -3.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.67'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '506.63'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.57'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.742'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.06'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000144'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.35%  '
$ws.Range('E27').Value = '  -3.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.47'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.42%  '
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.53'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.99'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.96'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.70'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.08%  '
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.05'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.11'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('E38').Value = '  +2.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '479.30'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.336'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '49.87'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.05'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.99'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.36'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.52'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.06%  '
$ws.Range('D46').Value = '2.920.92'
$ws.Range('E46').Value = '  -3.05%  '
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '140.19'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.84'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.42'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.00%  '
